$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Replace the formula in B1 (=160*10) with a plain static value
$ws.Range("B1").Value = 1000

# Move the active cell selection from B3 to B2
$ws.Range("B2").Select() | Out-Null
